$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Circumference") values added for the updated circumference algorithm
$ws.Range("C2").Value = 1119.134118676186
$ws.Range("C3").Value = 914.9747383594513
$ws.Range("C4").Value = 506.6589421033859
$ws.Range("C5").Value = 288.2497808933258
$ws.Range("C6").Value = 724.7493426799774
$ws.Range("C7").Value = 406.7594473361969
$ws.Range("C8").Value = 845.0752435922623
$ws.Range("C9").Value = 701.3523740768433
$ws.Range("C10").Value = 1160.721993088722
$ws.Range("C11").Value = 1213.567665100098
$ws.Range("C12").Value = 1619.498685359955
$ws.Range("C13").Value = 703.6122596263885
$ws.Range("C14").Value = 286.1492756605148
$ws.Range("C15").Value = 173.3797236680984
$ws.Range("C16").Value = 195.5807341337204
$ws.Range("C17").Value = 1430.805253386497
$ws.Range("C18").Value = 452.558436870575
$ws.Range("C19").Value = 216.3502861261368
$ws.Range("C20").Value = 704.4234417676926
$ws.Range("C21").Value = 242.2497808933258
$ws.Range("C22").Value = 650.7665876150131
$ws.Range("C23").Value = 1221.248904466629
$ws.Range("C24").Value = 306.433546423912
$ws.Range("C25").Value = 341.2792184352875
$ws.Range("C26").Value = 282.3919162750244
$ws.Range("C27").Value = 291.9482651948929
$ws.Range("C28").Value = 164.4091612100601
$ws.Range("C29").Value = 392.0315254926682
$ws.Range("C30").Value = 295.3208485841751
$ws.Range("C31").Value = 206.1076455116272
$ws.Range("C32").Value = 116.5685415267944
$ws.Range("C33").Value = 2429.586121559143
$ws.Range("C34").Value = 100.5685415267944
$ws.Range("C35").Value = 435.2447285652161
$ws.Range("C36").Value = 439.7300097942352
$ws.Range("C37").Value = 110.811182141304
$ws.Range("C38").Value = 267.3208485841751
$ws.Range("C39").Value = 60.76955199241638
$ws.Range("C40").Value = 362.7178171873093
$ws.Range("C41").Value = 936.3717069625854
$ws.Range("C42").Value = 395.7888848781586
$ws.Range("C43").Value = 305.3624787330627
$ws.Range("C44").Value = 365.64674949646
$ws.Range("C45").Value = 36.38477599620819
$ws.Range("C46").Value = 254.534051656723
$ws.Range("C47").Value = 418.1736608743668
$ws.Range("C48").Value = 295.3624787330627
$ws.Range("C49").Value = 1075.692555546761
$ws.Range("C50").Value = 496.6000670194626
$ws.Range("C51").Value = 648.3990565538406
$ws.Range("C52").Value = 238.4924215078354
$ws.Range("C53").Value = 675.5950146913528
$ws.Range("C54").Value = 407.7300097942352
$ws.Range("C55").Value = 627.068103313446
$ws.Range("C56").Value = 376.6589421033859
$ws.Range("C57").Value = 257.2203433513641
$ws.Range("C58").Value = 468.558436870575
$ws.Range("C59").Value = 392.2741661071777
$ws.Range("C60").Value = 243.5634891986847
$ws.Range("C61").Value = 230.4507913589478
$ws.Range("C62").Value = 722.4234417676926
$ws.Range("C63").Value = 216.2081507444382
$ws.Range("C64").Value = 1147.325024485588
$ws.Range("C65").Value = 1266.219466924667
$ws.Range("C66").Value = 397.7888848781586
$ws.Range("C67").Value = 187.9655101299286
$ws.Range("C68").Value = 1392.101716756821
$ws.Range("C69").Value = 302.4924215078354
$ws.Range("C70").Value = 819.6193999052048
$ws.Range("C71").Value = 289.9066350460052
$ws.Range("C72").Value = 319.2619735002518
$ws.Range("C73").Value = 249.1787132024765
$ws.Range("C74").Value = 551.9554054737091
$ws.Range("C75").Value = 234.1076455116272
$ws.Range("C76").Value = 357.4041088819504
$ws.Range("C77").Value = 994.8986183404922
$ws.Range("C78").Value = 344.7766922712326
$ws.Range("C79").Value = 499.6711347103119
$ws.Range("C80").Value = 349.2030984163284
$ws.Range("C81").Value = 327.2619735002518
$ws.Range("C82").Value = 466.1736608743668
$ws.Range("C83").Value = 387.8650048971176
$ws.Range("C84").Value = 624.4406867027283
$ws.Range("C85").Value = 822.4478269815445
$ws.Range("C86").Value = 44.62741661071777
$ws.Range("C87").Value = 238.6934319734573
$ws.Range("C88").Value = 372.1736608743668
$ws.Range("C89").Value = 312.0487704277039
$ws.Range("C90").Value = 235.3797236680984
$ws.Range("C91").Value = 280.8771975040436
$ws.Range("C92").Value = 257.4213538169861
$ws.Range("C93").Value = 207.82337474823
$ws.Range("C94").Value = 123.6396092176437
